## Adds two new "Variance" / "Correlation" sections (heading + R source
## code chunk + R console output) right before the existing
## "Données centrées réduites" Heading 1, mirroring the pre-existing
## "Covariance" section that sits just above the insertion point.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Inserts a brand-new, empty paragraph immediately before the paragraph
# currently living at index $idx, gives it style $pStyleName, and
# returns that new paragraph. Because InsertParagraphBefore re-seats the
# paragraph object at $idx onto the freshly created (empty) paragraph,
# the original content is pushed down to $idx + 1 — callers must bump
# their own tracking index after calling this.
function Insert-ParaBefore($idx, $pStyleName) {
    $doc = $word.ActiveDocument
    $anchorPara = $doc.Paragraphs($idx)
    $rng = $anchorPara.Range
    $rng.Collapse(1)
    $rng.InsertParagraphBefore()
    $newPara = $doc.Paragraphs($idx)
    $newPara.Range.Style = $pStyleName
    return $newPara
}

# Fills paragraph $para with a sequence of runs. $runs is an array of
# hashtables; each is either @{Text="..."; Style="..."} for a normal
# run (rendered with that character style) or @{Break=$true} for a
# <w:br/> line break.
function Add-Runs($para, $runs) {
    $doc = $word.ActiveDocument
    $rng = $para.Range
    $rng.Collapse(1)
    $start = $rng.Start

    $fullText = ""
    foreach ($r in $runs) {
        if ($r.Break -eq $true) {
            $fullText += [char]11
        } else {
            $fullText += $r.Text
        }
    }
    $rng.InsertAfter($fullText)

    $pos = $start
    foreach ($r in $runs) {
        if ($r.Break -eq $true) {
            $pos += 1
        } else {
            $len = $r.Text.Length
            if ($len -gt 0 -and $r.Style) {
                $subRng = $doc.Range($pos, $pos + $len)
                $subRng.Style = $r.Style
            }
            $pos += $len
        }
    }
}

# Wraps the (already populated) paragraph $para in a bookmark named
# $name, spanning its text (not the trailing paragraph mark).
function Add-ParaBookmark($para, $name) {
    $doc = $word.ActiveDocument
    $bmRng = $para.Range.Duplicate
    $bmRng.MoveEnd(1, -1)
    $doc.Bookmarks.Add($name, $bmRng)
}

# ---------------------------------------------------------------------
# Locate the insertion point: the "Données centrées réduites" Heading 1.
# ---------------------------------------------------------------------

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.StartsWith("Donn")) {
        $anchorIndex = $i
        break
    }
}

# ---------------------------------------------------------------------
# 1) "Variance (pour uniquement 2 colonnes)" heading
# ---------------------------------------------------------------------

$pVarHeading = Insert-ParaBefore $anchorIndex "Heading 2"
$anchorIndex += 1
Add-Runs $pVarHeading @(
    @{Text = "Variance (pour uniquement 2 colonnes)" }
)
Add-ParaBookmark $pVarHeading "variance-pour-uniquement-2-colonnes"

# ---------------------------------------------------------------------
# 2) var(x_matrix[,1:2]); source code line
# ---------------------------------------------------------------------

$pVarCode = Insert-ParaBefore $anchorIndex "Source Code"
$anchorIndex += 1
Add-Runs $pVarCode @(
    @{Text = "var"; Style = "KeywordTok" },
    @{Text = "(x_matrix[,"; Style = "NormalTok" },
    @{Text = "1"; Style = "DecValTok" },
    @{Text = ":"; Style = "OperatorTok" },
    @{Text = "2"; Style = "DecValTok" },
    @{Text = "]);"; Style = "NormalTok" }
)

# ---------------------------------------------------------------------
# 3) Console output of var(...)
# ---------------------------------------------------------------------

$pVarOut = Insert-ParaBefore $anchorIndex "Source Code"
$anchorIndex += 1
Add-Runs $pVarOut @(
    @{Text = "##                 X971.Guadeloupe X972.Martinique"; Style = "VerbatimChar" },
    @{Break = $true },
    @{Text = "## X971.Guadeloupe       185185011       147139734"; Style = "VerbatimChar" },
    @{Break = $true },
    @{Text = "## X972.Martinique       147139734       127959409"; Style = "VerbatimChar" }
)

# ---------------------------------------------------------------------
# 4) "Correlation (pour uniquement 2 colonnes)" heading
# ---------------------------------------------------------------------

$pCorHeading = Insert-ParaBefore $anchorIndex "Heading 2"
$anchorIndex += 1
Add-Runs $pCorHeading @(
    @{Text = "Correlation (pour uniquement 2 colonnes)" }
)
Add-ParaBookmark $pCorHeading "correlation-pour-uniquement-2-colonnes"

# ---------------------------------------------------------------------
# 5) cor(x_matrix[,1:2]) source code line
# ---------------------------------------------------------------------

$pCorCode = Insert-ParaBefore $anchorIndex "Source Code"
$anchorIndex += 1
Add-Runs $pCorCode @(
    @{Text = "cor"; Style = "KeywordTok" },
    @{Text = "(x_matrix[,"; Style = "NormalTok" },
    @{Text = "1"; Style = "DecValTok" },
    @{Text = ":"; Style = "OperatorTok" },
    @{Text = "2"; Style = "DecValTok" },
    @{Text = "])"; Style = "NormalTok" }
)

# ---------------------------------------------------------------------
# 6) Console output of cor(...)
# ---------------------------------------------------------------------

$pCorOut = Insert-ParaBefore $anchorIndex "Source Code"
$anchorIndex += 1
Add-Runs $pCorOut @(
    @{Text = "##                 X971.Guadeloupe X972.Martinique"; Style = "VerbatimChar" },
    @{Break = $true },
    @{Text = "## X971.Guadeloupe       1.0000000       0.9558526"; Style = "VerbatimChar" },
    @{Break = $true },
    @{Text = "## X972.Martinique       0.9558526       1.0000000"; Style = "VerbatimChar" }
)

# ---------------------------------------------------------------------
# Renumber the picture descr / chunk references: two new R chunks were
# inserted ahead of the existing ones, so "unnamed-chunk-9" -> -11,
# "unnamed-chunk-11" -> -13, "unnamed-chunk-12" -> -14 (processed in
# document order so the intermediate "-11" does not collide).
# ---------------------------------------------------------------------

$d.Content.Find.Execute("unnamed-chunk-9-1.png", $false, $false, $false, $false, $false, $true, 1, $false, "unnamed-chunk-11-1.png", 2) | Out-Null
$d.Content.Find.Execute("unnamed-chunk-11-1.png", $false, $false, $false, $false, $false, $true, 1, $false, "unnamed-chunk-13-1.png", 2) | Out-Null

$rng2 = $d.Range(0, 0)
$found = $rng2.Find.Execute("unnamed-chunk-13-1.png", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.End = $d.Content.End
$rng2.Find.Execute("unnamed-chunk-12-1.png", $false, $false, $false, $false, $false, $true, 1, $false, "unnamed-chunk-14-1.png", 2) | Out-Null
